$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.09
$summary.Range("B6").Value = 17
$summary.Range("B9").Value = 29.41

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 17
$status.Range("G4").Value = 29.41

# --- New trade row (#17) data shared by "All Trades" and "MarketMaking" sheets ---
$tradeNum     = 17
$tradeDate    = "2026-02-17"
$tradeTime    = "07:59:14"
$tradeStrat   = "MarketMaking"
$tradeSide    = "DOWN"
$tradeEntry   = 0.97
$tradeExit    = 0.97
$tradeStatus  = "CLOSED"
$tradePnlPct  = 0
$tradePnlUsd  = 0
$tradeCapital = 99.92
$tradeEntrySl = 0
$tradeExitSl  = 0
$tradeConf    = 0.6
$tradeEntryR  = "Normal spread capture: 19600 bps"
$tradeExitR   = "early_exit"
$tradeDur     = 0.13

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 18
    $ws.Cells.Item($row, 1).Value  = $tradeNum
    # Dates/times in this sheet are stored as plain text, not Excel date
    # serials, so force the cell to text before assigning (Excel would
    # otherwise auto-convert a recognizable date string into a date value).
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $tradeDate
    $dateCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value  = $tradeTime
    $ws.Cells.Item($row, 4).Value  = $tradeStrat
    $ws.Cells.Item($row, 5).Value  = $tradeSide
    $ws.Cells.Item($row, 6).Value  = $tradeEntry
    $ws.Cells.Item($row, 7).Value  = $tradeExit
    $ws.Cells.Item($row, 8).Value  = $tradeStatus
    $ws.Cells.Item($row, 9).Value  = $tradePnlPct
    $ws.Cells.Item($row, 10).Value = $tradePnlUsd
    $ws.Cells.Item($row, 11).Value = $tradeCapital
    $ws.Cells.Item($row, 12).Value = $tradeEntrySl
    $ws.Cells.Item($row, 13).Value = $tradeExitSl
    $ws.Cells.Item($row, 14).Value = $tradeConf
    $ws.Cells.Item($row, 15).Value = $tradeEntryR
    $ws.Cells.Item($row, 16).Value = $tradeExitR
    $ws.Cells.Item($row, 17).Value = $tradeDur
}
